# Plant data update: add ramp_rate column, drop the old "Wind" plants
# (154-157) and re-home the Hydro fleet from area 3 to area 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the four obsolete "Wind" plant rows (plant_id 154,155,156,157).
#    Everything below shifts up, so what was row 16 (plant_id 49900)
#    becomes row 12, etc.
$ws.Rows("12:15").Delete()

# 2. The remaining plants (now rows 12:36) are the Hydro fleet; move them
#    from area 3 to area 2.
for ($r = 12; $r -le 36; $r++) {
    $ws.Cells.Item($r, 2).Value2 = 2
}

# 3. Add the new "ramp_rate" column (G) with its header.
$ws.Range("G1").Value2 = "ramp_rate"

# Thermal plants (rows 2-11): ramp_rate values.
$ws.Cells.Item(2, 7).Value2 = 3.0
$ws.Cells.Item(3, 7).Value2 = 3.0
$ws.Cells.Item(4, 7).Value2 = 2.0
$ws.Cells.Item(5, 7).Value2 = 2.0
$ws.Cells.Item(6, 7).Value2 = 3.0
$ws.Cells.Item(7, 7).Value2 = 3.0
$ws.Cells.Item(8, 7).Value2 = 2.0
$ws.Cells.Item(9, 7).Value2 = 2.0
$ws.Cells.Item(10, 7).Value2 = 4.14
$ws.Cells.Item(11, 7).Value2 = 3.7

# Hydro plants (rows 12-36): ramp_rate is 100.0 for all of them.
for ($r = 12; $r -le 36; $r++) {
    $ws.Cells.Item($r, 7).Value2 = 100.0
}
